# Updated care data (v28)
#
# Applies the following content changes to the "routekaart status care" sheet:
#   1. Carint-Reggeland Groep (Stichting): status Voorlopig -> Vastgesteld
#   2. Remove the row for "Reinalda (Stichting Zorggroep, nu onder Kennemerhart)"
#   3. Add a new row for "Emergis (Stichting)" (status Vastgesteld), inserted
#      in alphabetical order (just before "Espria ...")
#   4. Add a new row for "Waardeburgh (Stichting)" (status Voorlopig), inserted
#      in alphabetical order (just before "Warande (Stichting)")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colA = $ws.Columns("A:A")

# 1. Update status for Carint-Reggeland Groep (Stichting)
$carint = $colA.Find("Carint-Reggeland Groep (Stichting)")
if ($carint) {
    $ws.Cells.Item($carint.Row, 2).Value = "Vastgesteld"
}

# 2. Remove the Reinalda row entirely (shifts subsequent rows up by one)
$reinalda = $colA.Find("Reinalda (Stichting Zorggroep, nu onder Kennemerhart)")
if ($reinalda) {
    $ws.Rows($reinalda.Row).Delete()
}

# 3. Insert Emergis (Stichting) right before Espria (...) keeping alphabetical order
$espria = $colA.Find("Espria (stichting ...) (met onderdelen Trans, Meander, Evean, GGZ Drenthe, icare)")
if ($espria) {
    $ws.Rows($espria.Row).Insert()
    $ws.Cells.Item($espria.Row, 1).Value = "Emergis (Stichting)"
    $ws.Cells.Item($espria.Row, 2).Value = "Vastgesteld"
}

# 4. Insert Waardeburgh (Stichting) right before Warande (Stichting) keeping alphabetical order
$warande = $colA.Find("Warande (Stichting)")
if ($warande) {
    $ws.Rows($warande.Row).Insert()
    $ws.Cells.Item($warande.Row, 1).Value = "Waardeburgh (Stichting)"
    $ws.Cells.Item($warande.Row, 2).Value = "Voorlopig"
}
